# Activities_Changes - 6th July 2023
# Update the Users sheet: replace the user "Nicole Bicho" with "Drew Koecher".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A2").Value = "Drew Koecher"
